# "Generate Report for Archive"
#
# The localization status report moved on from the handoff stage: the
# "Ready for handoff" status is now "In Translation" for every tracked
# file. Update the status cells on all three sheets and shrink the
# status columns so they still fit their (now shorter) contents, the
# same way Excel would after the text shrinks.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is mirrored in columns E (zh-cn) and F (de-de) ---
foreach ($cell in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    if ($wsOverview.Range($cell).Value2 -eq $oldStatus) {
        $wsOverview.Range($cell).Value = $newStatus
    }
}

# --- zh-cn / de-de sheets: status lives in column C ("Status") ---
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($cell in @("C2", "C3", "C4")) {
        if ($ws.Range($cell).Value2 -eq $oldStatus) {
            $ws.Range($cell).Value = $newStatus
        }
    }
}

# Re-fit the status columns now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
